$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "what is captured in this" descriptions for Woody Plant, Ground Cover Data, and Forest Management rows
$ws.Range("E4").Value = "Woody plants info - A more detail explanation is provided in the Understanding data - Edem markdown document"
$ws.Range("E6").Value = "Ground cover plants info - A more detail explanation is provided in the Understanding data - Edem markdown document"
$ws.Range("E8").Value = "Forest management effort info -  A more detail explanation is provided in the Understanding data - Edem markdown document"

# Enable wrap text for these cells and set row heights
$ws.Range("E4").WrapText = $true
$ws.Range("E6").WrapText = $true
$ws.Range("E8").WrapText = $true

$ws.Rows.Item(4).RowHeight = 45
$ws.Rows.Item(6).RowHeight = 45
$ws.Rows.Item(8).RowHeight = 45

$ws.Range("G14").Select()
